$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in row 30: 24.11.18, 10:00 -> 16:00, "-squirrel"
$ws.Range("A30").Value = "24.11.18"
$ws.Range("B30").Value = 0.41666666666666669
$ws.Range("C30").Value = 0.66666666666666663
$ws.Range("E30").Value = "-squirrel"

# Fill in row 31: 26.11.18, 11:00 -> 18:00
$ws.Range("A31").Value = "26.11.18"
$ws.Range("B31").Value = 0.45833333333333331
$ws.Range("C31").Value = 0.75
$ws.Range("E31").Value = "-Output Options"

# Extend the D column elapsed-time formula down through the new rows
# (mirrors dragging the fill handle from D29 to D31, continuing the
# existing shared formula C#-B#)
$ws.Range("D30:D31").Formula = "=C30-B30"

# Update selection to A32
$ws.Range("A32").Select()
